$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (row 2 through row 10) from serial date 45179 to 45180
# (i.e. bump the "Förändrad" date forward by one day for each entry)
$ws.Range("C2:C10").Value = 45180
